$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 5266105.5
$ws.Range("I137").Value = 10002240
$ws.Range("J137").Value = 3733.6667
$ws.Range("K137").Value = 30006720
$ws.Range("L137").Value = 11201.0001
$ws.Range("M137").Value = -30004170
$ws.Range("N137").Value = -16301.0001
$ws.Range("H138").Value = 1919414
$ws.Range("I138").Value = 2294.0417
$ws.Range("J138").Value = 2649745.2
$ws.Range("K138").Value = 6882.125100000001
$ws.Range("L138").Value = 7949235.600000001
$ws.Range("N138").Value = -7959515.600000001
$ws.Range("M138").Value = -1742.125100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4996304
$ws.Range("I32").Value = 5767395
$ws.Range("J32").Value = 27051
$ws.Range("K32").Value = 5767395
$ws.Range("L32").Value = 27051
$ws.Range("M32").Value = -5767108
$ws.Range("N32").Value = -27625
$ws.Range("H45").Value = 2784.4348
$ws.Range("I45").Value = 2897.4211
$ws.Range("J45").Value = 2247.75
$ws.Range("K45").Value = 2897.4211
$ws.Range("L45").Value = 2247.75
$ws.Range("M45").Value = -2520.4211
$ws.Range("N45").Value = -3001.75
$ws.Range("H61").Value = 71572680
$ws.Range("I61").Value = 111223500
$ws.Range("J61").Value = 201206.8
$ws.Range("K61").Value = 111223500
$ws.Range("L61").Value = 201206.8
$ws.Range("M61").Value = -111223288
$ws.Range("N61").Value = -201630.8
$ws.Range("H74").Value = 10501186
$ws.Range("I74").Value = 14765618
$ws.Range("K74").Value = 14765618
$ws.Range("M74").Value = -14764744
$ws.Range("H77").Value = 10501186
$ws.Range("I77").Value = 14765618
$ws.Range("K77").Value = 73828090
$ws.Range("M77").Value = -73823722
$ws.Range("H122").Value = 6175321.5
$ws.Range("I122").Value = 2713.1538
$ws.Range("K122").Value = 8139.4614
$ws.Range("M122").Value = -5689.4614
$ws.Range("H136").Value = 71572680
$ws.Range("I136").Value = 111223500
$ws.Range("J136").Value = 201206.8
$ws.Range("K136").Value = 333670500
$ws.Range("L136").Value = 603620.3999999999
$ws.Range("M136").Value = -333667950
$ws.Range("N136").Value = -608720.3999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3249.275
$ws.Range("I134").Value = 2726.5757
$ws.Range("J134").Value = 5713.4287
$ws.Range("K134").Value = 8179.7271
$ws.Range("L134").Value = 17140.2861
$ws.Range("M134").Value = -5644.7271
$ws.Range("N134").Value = -22210.2861

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 7800
$ws.Range("J62").Value = 8000
$ws.Range("L62").Value = 8000
$ws.Range("N62").Value = -9248
$ws.Range("H65").Value = 7800
$ws.Range("J65").Value = 8000
$ws.Range("L65").Value = 40000
$ws.Range("N65").Value = -46240
$ws.Range("H132").Value = 58673.055
$ws.Range("I132").Value = 2723.7856
$ws.Range("J132").Value = 254495.5
$ws.Range("K132").Value = 8171.3568
$ws.Range("L132").Value = 763486.5
$ws.Range("M132").Value = -5641.3568
$ws.Range("N132").Value = -768546.5
$ws.Range("H134").Value = 29398.256
$ws.Range("I134").Value = 1672.72
$ws.Range("J134").Value = 78908.14
$ws.Range("K134").Value = 5018.16
$ws.Range("L134").Value = 236724.42
$ws.Range("M134").Value = -2483.16
$ws.Range("N134").Value = -241794.42

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1997.75
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 1997.75
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 5993.25
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -7615.25
$ws.Range("H71").Value = 1997.75
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 1997.75
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 17979.75
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -26091.75
$ws.Range("H120").Value = 5790
$ws.Range("I120").Value = 5790
$ws.Range("J120").Value = 0
$ws.Range("K120").Value = 17370
$ws.Range("L120").Value = 0
$ws.Range("M120").Value = -12532
$ws.Range("N120").ClearContents()
$ws.Range("H131").Value = 790.9355
$ws.Range("J131").Value = 922.2273
$ws.Range("L131").Value = 2766.6819
$ws.Range("N131").Value = -12846.6819

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 14996.667
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 14996.667
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 14996.667
$ws.Range("M57").ClearContents()
$ws.Range("N57").Value = -16636.667
$ws.Range("H102").Value = 2076.1052
$ws.Range("I102").Value = 2146.7273
$ws.Range("J102").Value = 1979
$ws.Range("K102").Value = 2146.7273
$ws.Range("L102").Value = 1979
$ws.Range("M102").Value = -524.7273
$ws.Range("N102").Value = -5223
$ws.Range("H109").Value = 11190
$ws.Range("J109").Value = 11190
$ws.Range("L109").Value = 11190
$ws.Range("N109").Value = -13270
$ws.Range("H122").Value = 3552.5
$ws.Range("I122").Value = 2938.5
$ws.Range("J122").Value = 4166.5
$ws.Range("K122").Value = 8815.5
$ws.Range("L122").Value = 12499.5
$ws.Range("M122").Value = -6365.5
$ws.Range("N122").Value = -17399.5
$ws.Range("H124").Value = 50780
$ws.Range("J124").Value = 50780
$ws.Range("L124").Value = 50780
$ws.Range("N124").Value = -60600
$ws.Range("H132").Value = 61652.824
$ws.Range("I132").Value = 52590.2
$ws.Range("J132").Value = 74599.42999999999
$ws.Range("K132").Value = 157770.6
$ws.Range("L132").Value = 223798.29
$ws.Range("M132").Value = -155240.6
$ws.Range("N132").Value = -228858.29
$ws.Range("H135").Value = 57645
$ws.Range("J135").Value = 57645
$ws.Range("L135").Value = 57645
$ws.Range("N135").Value = -67785

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3563.75
$ws.Range("I122").Value = 3351.6667
$ws.Range("J122").Value = 4200
$ws.Range("K122").Value = 10055.0001
$ws.Range("L122").Value = 12600
$ws.Range("M122").Value = -7605.000100000001
$ws.Range("N122").Value = -17500
$ws.Range("H132").Value = 40036.184
$ws.Range("I132").Value = 1956.6923
$ws.Range("J132").Value = 75395.71000000001
$ws.Range("K132").Value = 5870.0769
$ws.Range("L132").Value = 226187.13
$ws.Range("M132").Value = -3340.0769
$ws.Range("N132").Value = -231247.13
$ws.Range("H136").Value = 71451.78999999999
$ws.Range("I136").Value = 34347.582
$ws.Range("J136").Value = 176018.19
$ws.Range("K136").Value = 103042.746
$ws.Range("L136").Value = 528054.5700000001
$ws.Range("M136").Value = -100492.746
$ws.Range("N136").Value = -533154.5700000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2749.9524
$ws.Range("I122").Value = 2041.9231
$ws.Range("J122").Value = 3900.5
$ws.Range("K122").Value = 6125.7693
$ws.Range("L122").Value = 11701.5
$ws.Range("M122").Value = -3675.7693
$ws.Range("N122").Value = -16601.5
$ws.Range("H132").Value = 48409.977
$ws.Range("I132").Value = 35028.4
$ws.Range("J132").Value = 79290.53999999999
$ws.Range("K132").Value = 105085.2
$ws.Range("L132").Value = 237871.62
$ws.Range("M132").Value = -102555.2
$ws.Range("N132").Value = -242931.62
$ws.Range("H136").Value = 53182.82
$ws.Range("I136").Value = 31270.303
$ws.Range("K136").Value = 93810.909
$ws.Range("M136").Value = -91260.909
